$d = $word.ActiveDocument

# 1. Change tense: "Leverages" -> "Leveraged" (same length, keeps offsets stable)
$d.Content.Find.Execute("Leverages", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Leveraged", 2) | Out-Null

# 2. Re-find "Leveraged" to get its exact range (end position is where the
#    _GoBack bookmark should now live, splitting the run in two).
$rng = $d.Content
$rng.Find.Execute("Leveraged", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0) | Out-Null

# 3. Move the _GoBack bookmark from the end of the previous paragraph to
#    right after "Leveraged".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bookmarkRange = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
